# "Circulation Course - ajout et fin"
#  - Fin de la page circulation en course
#  - Traduction en anglais
#  - Ajout des images tirees du guide d'organisation UCI

$wb = $excel.ActiveWorkbook

# --- COMM sheet: finalize the Commissaires panel chair entry -------------
# "Presidente" -> "President" (role title correction) and the short
# placeholder name "Thierry D" is completed to "Thierry Diederen".
$comm = $wb.Worksheets.Item("COMM")
$comm.Range("B2").Value = "Président du Collège des commissaires "
$comm.Range("C2").Value = "Thierry Diederen"

# --- SOUTIEN sheet: tighten up row heights now that the page is done -----
$soutien = $wb.Worksheets.Item("SOUTIEN")
$soutien.Rows.Item(5).RowHeight = 17
$soutien.Rows.Item(9).RowHeight = 17

# --- Leave the workbook focused on the COMM sheet, selection on B11 ------
$comm.Activate()
$comm.Range("B11").Select()
